# Auto-update draw results: append the 2025-12-18 Pick 3 result as a new
# row (93) at the bottom of the results table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# New draw data for 2025-12-18.
$date        = "2025-12-18"
$game        = "Pick 3"
$phase       = "251218"
$result      = "6-5-4"
$insertedAt  = "2025-12-18T21:45:12.708+04:00"

# Columns A (date) and C (phase) look like a date / a plain number to the
# recalculation engine, so force them to Text before entry (matching the
# rest of the sheet, which stores every value as literal text) and then
# clear the formatting again so no new cell style gets attached - exactly
# like the existing A1:E92 cells, which carry no explicit style index.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = $date
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = $game

$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = $phase
$ws.Range("C$row").ClearFormats()

$ws.Range("D$row").Value = $result

$ws.Range("E$row").Value = $insertedAt
